{"js": "// Update the invoice: shift the document date by one day, rename the\n// item, and correct the price / total / amount-in-words to match.\n\nasync function replaceAll(context, searchText, replacement, options) {\n  const results = context.document.body.search(searchText, options || { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n  return results.items.length;\n}\n\n// 1) Date \"26.5.2024\" -> \"27.5.2024\" (appears twice in the document).\nawait replaceAll(context, \"26.5.2024\", \"27.5.2024\");\n\n// 2) Item designation \"YH9000AE\" -> \"KALTMANN MASCHINEN 2\".\nawait replaceAll(context, \"YH9000AE\", \"KALTMANN MASCHINEN 2\");\n\n// 3) Price / sum \"18070.08\" -> \"14450.00\" (appears three times).\nawait replaceAll(context, \"18070.08\", \"14450.00\");\n\n// 4) Total amount spelled out in words (Ukrainian).\nawait replaceAll(\n  context,\n  \"\u0432\u0456\u0441\u0456\u043c\u043d\u0430\u0434\u0446\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0441\u0435\u043c\u044c\u0434\u0435\u0441\u044f\u0442 \u0433\u0440\u0438\u0432\u0435\u043d\u044c \u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\",\n  \"\u0447\u043e\u0442\u0438\u0440\u043d\u0430\u0434\u0446\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0447\u043e\u0442\u0438\u0440\u0438\u0441\u0442\u0430 \u043f\u044f\u0442\u044c\u0434\u0435\u0441\u044f\u0442 \u0433\u0440\u0438\u0432\u0435\u043d\u044c \u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\"\n);\n", "ps1": "# Update the invoice: shift the document date by one day, rename the\n# item, and correct the price / total / amount-in-words to match.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n# 1) Date \"26.5.2024\" -> \"27.5.2024\" (appears twice in the document).\n$d.Content.Find.Execute(\"26.5.2024\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"27.5.2024\", $wdReplaceAll)\n\n# 2) Item designation \"YH9000AE\" -> \"KALTMANN MASCHINEN 2\".\n$d.Content.Find.Execute(\"YH9000AE\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"KALTMANN MASCHINEN 2\", $wdReplaceAll)\n\n# 3) Price / sum \"18070.08\" -> \"14450.00\" (appears three times).\n$d.Content.Find.Execute(\"18070.08\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"14450.00\", $wdReplaceAll)\n\n# 4) Total amount spelled out in words (Ukrainian).\n$d.Content.Find.Execute(\"\u0432\u0456\u0441\u0456\u043c\u043d\u0430\u0434\u0446\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0441\u0435\u043c\u044c\u0434\u0435\u0441\u044f\u0442 \u0433\u0440\u0438\u0432\u0435\u043d\u044c \u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"\u0447\u043e\u0442\u0438\u0440\u043d\u0430\u0434\u0446\u044f\u0442\u044c \u0442\u0438\u0441\u044f\u0447 \u0447\u043e\u0442\u0438\u0440\u0438\u0441\u0442\u0430 \u043f\u044f\u0442\u044c\u0434\u0435\u0441\u044f\u0442 \u0433\u0440\u0438\u0432\u0435\u043d\u044c \u043d\u0443\u043b\u044c \u043a\u043e\u043f\u0456\u0439\u043e\u043a\", $wdReplaceAll)\n"}
